$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'SM 画廊 自殺などに過労自殺をなくすためには 8月に起こりきましょうか'
$ws.Range("C3").Value = '宿 意見が意見が浮かんだ人からじゃあ僕からですか 過労自殺をなくすための対策としてはまず 企業が過労自殺をさせない環境を作るのが一番手っ取り早い方法だと思いますよ じゃあ どうやって過労自殺を起こさないような環境を作り上げるかと言うと あの国とかがあの 過労自殺をするくらいなら その会社をぶっつぶせ みたいな理論をあげて 過労自殺に追い込まれた若者を追い込んだ企業を報告する報告しちゃっていいよ みたいな制度を作って そしたら 企業 それを抑止するために'
$ws.Range("C4").Value = 'それを抑止するためにまともな労働環境を作らざるを得へんから なんかそうやったらできるんじゃないかな 過労自殺は防げるんじゃないかな 自分が追い込まれて死ぬくらいなら 企業のお人形さんみたいな そんな感じなんだこれからじゃあいつかの意見としては面白いから これは働いてるから起きるわけで そのことも 誰も当たらなくていいよ にるって言う なんでなんか最近 AI とかで人間の代わりとかできる 待ってるから'
$ws.Range("C5").Value = '音が鳴ってるから 極力その なんて言うかなことを 人間じゃなくては 機械とか家に任せてもやることによって負けるな クマ なくす お神輿 ないけども 減らすことはできるんじゃないか っていう人ですが 他に行ける人は国 コロナ 対策として 今きが違うけどそのような 角8 長時間労働のその取り組みを強化するとこがどこでも僕は対策に物言ってた民族んのアン'
$ws.Range("C7").Value = 'アロー歯科 だってマジで帽子をなくしてる人に対して めちゃめちゃ お金をあげるんだと その'
$ws.Range("C8").Value = '働いてる人が もう俺からしてないんじゃないみたいな感じに思わせるレベルのお金はあげる行きたいものは見合った金額をあげるって言う事からお金があるから辞めてもすぐ辞めても問題ないみたいな言い方'
$ws.Range("C9").Value = '後は 茅ヶ崎 お父さん やっぱの中 それぞれのなんか会 従業員って言った会社員が なんか 労働環境を見れる あれ なんかなんかしかみたいなんで来たらなくなるんじゃないかなと思って混んでて働いてますよ みたいな なんか 結構 働いてる人がいれば ちょっと声かけてなんかちょっとなんかあれかな とかなとかってさせるような企業がそんなマネジメントしてるとは思えへんねんな これからどうなるかわからん'
$ws.Range("C10").Value = 'ロベカル 動画から 何て言ったらさ ママは残り6分 何か他に意見がある人は 企業報告制度下で PAO 銚子 上がらない人のその後のことを考えると被害に遭ってない人に似合う形になるので それは現場の'
$ws.Range("C11").Value = 'それは現場の仲間意識に繋がる環境になってる人 見つけたら自分に飛び火するかもしれないって言う さまあ 期間と言うか 創傷感に駆られてさあの自分の身を守るための 偽りの正義 みたいな なんて言うの 自分に被害を被るのは嫌やからなんかその過労死しそうな人は仕方なく 守る みたいな なんか 報告したりとかして36とかと一緒でその社長からその言えない みたいな'
$ws.Range("C12").Value = '言えない みたいなことが起こる可能性はないことはないといけんって大きいを辞めた後のことを考えてるとこじゃないって考えてないけどお金もらっても 企業の方針 変わらんから よかろうしたんじゃないか みたいな感じでお金もらってるから 耐えられる'
$ws.Range("C13").Value = 'おまけ もらえる形がみたいな感じて自殺にはならんかなって言うか 前 あと3分版 5分版 担当 羽柴さんは解決にはならん つながる'
$ws.Range("C14").Value = 'ながら気はするな なんだったっけ 従業員の労働環境とか 時間を可視化するって言う それ プラスが追いつかなんかで可視化したやつを監視する人を助けるとかなんかそれが多くの人間がやるみたいな監査が 労働環境とかを観察するのはだいぶになった 金具 付けなきゃいけないみたいな 難しい'
$ws.Range("C15").Value = '後で考えるなら 企業側にどういう なんか リスクなんてやらその改善に努めます だけじゃあかんや そんな そういう気持ち かに それでひかかったとしてな 麻疹に辛そうなってしまった時の対応のことまで考えないと結構 緊張した難しそうなんかそれでなんかあの関さんに引っかかった企業はなんかかかるような法人税 増やす とか あーそれそれでもいたらもう'
$ws.Range("C16").Value = '誘わないないないない 傾向 納豆 から 第三者機関がそんな中で内部告発 れたやつか行かないとかなんかだとかそういう 縛りをつけて 勤怠管理と言うか'
$ws.Range("C17").Value = '勤怠管理というか'

# Clear C6 content (target keeps an empty inlineStr cell)
$ws.Range("C6").ClearContents()

# Remove rows 18:20 (old rows 18,19,20 data no longer present)
$ws.Rows("18:20").Delete()

